$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $style = $rng.Style
    $rng.Value = "'" + $val
    $rng.Style = $style
}

function Set-PlainText($ws, $cellAddr, $val) {
    $ws.Range($cellAddr).Value = $val
}

Set-TextValue $ws "D2" "315.32"
Set-TextValue $ws "E2" "2.60%"
Set-TextValue $ws "D3" "39.45"
Set-TextValue $ws "E3" "2.31%"
Set-TextValue $ws "D4" "5.140"
Set-TextValue $ws "E4" "0.75%"
Set-TextValue $ws "D5" "0.08182"
Set-TextValue $ws "E5" "0.87%"
Set-TextValue $ws "D6" "1.963"
Set-TextValue $ws "E6" "-0.03%"
Set-TextValue $ws "D7" "8.216"
Set-TextValue $ws "E7" "3.37%"
Set-TextValue $ws "D8" "0.9261"
Set-TextValue $ws "E8" "-0.34%"
Set-TextValue $ws "D9" "0.1403"
Set-TextValue $ws "E9" "-1.82%"
Set-TextValue $ws "D10" "0.1980"
Set-TextValue $ws "E10" "1.11%"
Set-TextValue $ws "D11" "0.09042"
Set-TextValue $ws "E11" "-0.38%"
Set-TextValue $ws "D12" "0.03510"
Set-TextValue $ws "E12" "0.05%"
Set-TextValue $ws "D13" "0.09827"
Set-TextValue $ws "E13" "0.00%"
Set-TextValue $ws "D14" "0.001391"
Set-TextValue $ws "E14" "-1.04%"
Set-TextValue $ws "D15" "0.005972"
Set-TextValue $ws "E15" "-3.19%"
Set-TextValue $ws "D17" "4.238"
Set-TextValue $ws "E17" "1.03%"
Set-TextValue $ws "D18" "3.137"
Set-TextValue $ws "E18" "-8.23%"
Set-TextValue $ws "D19" "0.3464"
Set-TextValue $ws "E19" "0.05%"
Set-TextValue $ws "D20" "0.1345"
Set-TextValue $ws "E20" "0.24%"
Set-TextValue $ws "D21" "4.764"
Set-TextValue $ws "E21" "-0.58%"
Set-TextValue $ws "D22" "0.2428"
Set-TextValue $ws "E22" "-1.07%"
Set-TextValue $ws "D23" "0.04383"
Set-TextValue $ws "E23" "0.48%"
Set-TextValue $ws "E24" "-0.06%"
Set-TextValue $ws "D25" "0.004787"
Set-TextValue $ws "E25" "-0.98%"
Set-TextValue $ws "E26" "-0.16%"
Set-TextValue $ws "D27" "0.0003998"
Set-TextValue $ws "E27" "-10.10%"
Set-TextValue $ws "D39" "0.02180"
Set-TextValue $ws "E39" "4.52%"
Set-TextValue $ws "D40" "0.05186"
Set-TextValue $ws "E40" "1.23%"
Set-TextValue $ws "D41" "0.007553"
Set-TextValue $ws "E41" "0.96%"
Set-TextValue $ws "D42" "0.009795"
Set-TextValue $ws "E42" "-3.34%"
Set-TextValue $ws "D43" "0.1375"
Set-TextValue $ws "E43" "1.42%"
Set-TextValue $ws "D44" "0.002129"
Set-TextValue $ws "E44" "-0.16%"
Set-TextValue $ws "D45" "0.009128"
Set-TextValue $ws "E45" "-1.67%"
Set-TextValue $ws "E46" "2.28%"
Set-TextValue $ws "D47" "0.00000000750"
Set-TextValue $ws "E47" "-0.17%"
Set-PlainText $ws "B48" "BOLO"
Set-PlainText $ws "C48" "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
Set-TextValue $ws "D48" "0.002763"
Set-TextValue $ws "E48" "-8.76%"
Set-PlainText $ws "B49" "CoinbaseStockToken"
Set-PlainText $ws "C49" "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
Set-TextValue $ws "D49" "0.001200"
Set-TextValue $ws "E49" "-25.08%"
Set-TextValue $ws "D50" "0.00002099"
Set-TextValue $ws "E50" "-0.17%"
Set-TextValue $ws "D51" "0.0001999"
Set-TextValue $ws "E51" "-0.17%"
